# Update Betfair Back/Lay odds values for 2025-12-25 games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.08
$ws.Range("O2").Value = 1.37
$ws.Range("Q2").Value = 2.1
$ws.Range("T2").Value = 1.85
$ws.Range("U2").Value = 1.97
$ws.Range("X2").Value = 14.5
$ws.Range("AB2").Value = 15.5

# Row 3
$ws.Range("L3").Value = 1.42
$ws.Range("P3").Value = 1.79
$ws.Range("Q3").Value = 1.94
$ws.Range("S3").Value = 3.65
$ws.Range("AB3").Value = 990

# Row 5
$ws.Range("L5").Value = 1.32
$ws.Range("T5").Value = 1.92

# Row 6
$ws.Range("G6").Value = 3.3
$ws.Range("H6").Value = 2.32
$ws.Range("I6").Value = 2.56
$ws.Range("J6").Value = 3.35
$ws.Range("N6").Value = 3.65
$ws.Range("P6").Value = 1.92
$ws.Range("Q6").Value = 1.88
$ws.Range("S6").Value = 3.25
$ws.Range("U6").Value = 2.14
$ws.Range("V6").Value = 1.64
$ws.Range("W6").Value = 1.43
$ws.Range("X6").Value = 18
$ws.Range("AH6").Value = 17.5
$ws.Range("AI6").Value = 980
$ws.Range("AK6").Value = 38
$ws.Range("AL6").Value = 980

# Row 7
$ws.Range("Q7").Value = 1.66
$ws.Range("S7").Value = 2.66
$ws.Range("U7").Value = 1.91
$ws.Range("AB7").Value = 13

# Row 8
$ws.Range("F8").Value = 1.64
$ws.Range("K8").Value = 4.2
$ws.Range("T8").Value = 2.32

# Row 9
$ws.Range("O9").Value = 1.01
$ws.Range("P9").Value = 1.36
